# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AB$rowA")
    $rangeB = $ws.Range("B$rowB`:AB$rowB")
    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()
    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# Rows 105 and 106 (match ids 6077763 / 6077497) had their data swapped
Swap-Rows $ws 105 106

# Rows 113 and 114 (match ids 6078263 / 6078996) had their data swapped
Swap-Rows $ws 113 114

# Rows 139 and 140 (match ids 7723531 / 7723532) had their data swapped
Swap-Rows $ws 139 140

# Row 228: updated closing odds
$ws.Range("M228").Value = 1.9
$ws.Range("O228").Value = 4
